$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 171.25
$ws.Range("I9").Value = 147.6
$ws.Range("J9").Value = 210.66667
$ws.Range("K9").Value = 147.6
$ws.Range("L9").Value = 210.66667
$ws.Range("M9").Value = 21.40000000000001
$ws.Range("N9").Value = -548.6666700000001
$ws.Range("H15").Value = 177307
$ws.Range("I15").Value = 177307
$ws.Range("K15").Value = 531921
$ws.Range("M15").Value = -531752
$ws.Range("H76").Value = 3706696.5
$ws.Range("I76").Value = 4276389
$ws.Range("J76").Value = 3695.75
$ws.Range("K76").Value = 4276389
$ws.Range("L76").Value = 3695.75
$ws.Range("M76").Value = -4276074
$ws.Range("N76").Value = -4325.75
$ws.Range("H79").Value = 3706696.5
$ws.Range("I79").Value = 4276389
$ws.Range("J79").Value = 3695.75
$ws.Range("K79").Value = 4276389
$ws.Range("L79").Value = 3695.75
$ws.Range("M79").Value = -4275297
$ws.Range("N79").Value = -5879.75
$ws.Range("H107").Value = 1389843.2
$ws.Range("I107").Value = 2778727.8
$ws.Range("J107").Value = 958.75
$ws.Range("K107").Value = 2778727.8
$ws.Range("L107").Value = 958.75
$ws.Range("M107").Value = -2776807.8
$ws.Range("N107").Value = -4798.75
$ws.Range("H109").Value = 16666.666
$ws.Range("J109").Value = 16666.666
$ws.Range("L109").Value = 16666.666
$ws.Range("N109").Value = -19440.666
$ws.Range("H132").Value = 364823.44
$ws.Range("I132").Value = 435064.78
$ws.Range("K132").Value = 1305194.34
$ws.Range("M132").Value = -1302664.34
$ws.Range("H138").Value = 3942410.5
$ws.Range("I138").Value = 2179417.2
$ws.Range("J138").Value = 4331643
$ws.Range("K138").Value = 6538251.600000001
$ws.Range("L138").Value = 12994929
$ws.Range("M138").Value = -6533111.600000001
$ws.Range("N138").Value = -13005209
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2354.5625
$ws.Range("I61").Value = 1322.909
$ws.Range("K61").Value = 1322.909
$ws.Range("M61").Value = -1110.909
$ws.Range("H63").Value = 23232.5
$ws.Range("I63").Value = 27279
$ws.Range("J63").Value = 3000
$ws.Range("K63").Value = 27279
$ws.Range("L63").Value = 3000
$ws.Range("M63").Value = -26593
$ws.Range("N63").Value = -4372
$ws.Range("H66").Value = 23232.5
$ws.Range("I66").Value = 27279
$ws.Range("J66").Value = 3000
$ws.Range("K66").Value = 136395
$ws.Range("L66").Value = 15000
$ws.Range("M66").Value = -132963
$ws.Range("N66").Value = -21864
$ws.Range("H102").Value = 1642
$ws.Range("I102").Value = 1642
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1642
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -20
$ws.Range("N102").ClearContents()
$ws.Range("H110").Value = 476.37036
$ws.Range("I110").Value = 454.86957
$ws.Range("J110").Value = 600
$ws.Range("K110").Value = 454.86957
$ws.Range("L110").Value = 600
$ws.Range("M110").Value = 1590.13043
$ws.Range("N110").Value = -4690
$ws.Range("H136").Value = 2354.5625
$ws.Range("I136").Value = 1322.909
$ws.Range("K136").Value = 3968.727
$ws.Range("M136").Value = -1418.727

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3085.0938
$ws.Range("I105").Value = 2884.2917
$ws.Range("J105").Value = 3687.5
$ws.Range("K105").Value = 2884.2917
$ws.Range("L105").Value = 3687.5
$ws.Range("M105").Value = -1137.2917
$ws.Range("N105").Value = -7181.5
$ws.Range("H107").Value = 969.4706
$ws.Range("I107").Value = 628
$ws.Range("J107").Value = 1457.2858
$ws.Range("K107").Value = 628
$ws.Range("L107").Value = 1457.2858
$ws.Range("M107").Value = 1292
$ws.Range("N107").Value = -5297.2858
$ws.Range("H134").Value = 2862.6667
$ws.Range("I134").Value = 2260.12
$ws.Range("K134").Value = 6780.36
$ws.Range("M134").Value = -4245.36

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1251.75
$ws.Range("I31").Value = 961.96
$ws.Range("J31").Value = 3666.6667
$ws.Range("K31").Value = 961.96
$ws.Range("L31").Value = 3666.6667
$ws.Range("M31").Value = -666.96
$ws.Range("N31").Value = -4256.6667
$ws.Range("H34").Value = 1251.75
$ws.Range("I34").Value = 961.96
$ws.Range("J34").Value = 3666.6667
$ws.Range("K34").Value = 961.96
$ws.Range("L34").Value = 3666.6667
$ws.Range("M34").Value = -759.96
$ws.Range("N34").Value = -4070.6667
$ws.Range("H62").Value = 35713.285
$ws.Range("I62").Value = 70000
$ws.Range("J62").Value = 9998.25
$ws.Range("K62").Value = 70000
$ws.Range("L62").Value = 9998.25
$ws.Range("M62").Value = -69376
$ws.Range("N62").Value = -11246.25
$ws.Range("H65").Value = 35713.285
$ws.Range("I65").Value = 70000
$ws.Range("J65").Value = 9998.25
$ws.Range("K65").Value = 350000
$ws.Range("L65").Value = 49991.25
$ws.Range("M65").Value = -346880
$ws.Range("N65").Value = -56231.25
$ws.Range("H99").Value = 15626488
$ws.Range("I99").Value = 62500000
$ws.Range("J99").Value = 1983.3334
$ws.Range("K99").Value = 62500000
$ws.Range("L99").Value = 1983.3334
$ws.Range("M99").Value = -62498502
$ws.Range("N99").Value = -4979.3334
$ws.Range("H105").Value = 548.9
$ws.Range("I105").Value = 537.5
$ws.Range("J105").Value = 594.5
$ws.Range("K105").Value = 537.5
$ws.Range("L105").Value = 594.5
$ws.Range("M105").Value = 1209.5
$ws.Range("N105").Value = -4088.5
$ws.Range("H126").Value = 15626488
$ws.Range("I126").Value = 62500000
$ws.Range("J126").Value = 1983.3334
$ws.Range("K126").Value = 187500000
$ws.Range("L126").Value = 5950.0002
$ws.Range("M126").Value = -187497530
$ws.Range("N126").Value = -10890.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("M98").ClearContents()
$ws.Range("N98").ClearContents()
$ws.Range("H114").Value = 1181
$ws.Range("J114").Value = 1678.75
$ws.Range("L114").Value = 5036.25
$ws.Range("N114").Value = -11544.25
$ws.Range("H117").Value = 837.6923
$ws.Range("I117").Value = 893
$ws.Range("J117").Value = 790.2857
$ws.Range("K117").Value = 2679
$ws.Range("L117").Value = 2370.8571
$ws.Range("M117").Value = 763
$ws.Range("N117").Value = -9254.857100000001
$ws.Range("H132").Value = 1183.8636
$ws.Range("I132").Value = 816.5
$ws.Range("J132").Value = 1393.7858
$ws.Range("K132").Value = 7348.5
$ws.Range("L132").Value = 12544.0722
$ws.Range("M132").Value = -4818.5
$ws.Range("N132").Value = -17604.0722
$ws.Range("H136").Value = 1620.2
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2615.1428
$ws.Range("I80").Value = 2700
$ws.Range("J80").Value = 2601
$ws.Range("K80").Value = 2700
$ws.Range("L80").Value = 2601
$ws.Range("M80").Value = -1702
$ws.Range("N80").Value = -4597
$ws.Range("H83").Value = 2615.1428
$ws.Range("I83").Value = 2700
$ws.Range("J83").Value = 2601
$ws.Range("K83").Value = 13500
$ws.Range("L83").Value = 13005
$ws.Range("M83").Value = -8508
$ws.Range("N83").Value = -22989
$ws.Range("H113").Value = 1533.909
$ws.Range("I113").Value = 1566.6666
$ws.Range("J113").Value = 1494.6
$ws.Range("K113").Value = 1566.6666
$ws.Range("L113").Value = 1494.6
$ws.Range("M113").Value = 603.3334
$ws.Range("N113").Value = -5834.6
$ws.Range("H122").Value = 2223921.5
$ws.Range("I122").Value = 5556554
$ws.Range("J122").Value = 2166.6667
$ws.Range("K122").Value = 16669662
$ws.Range("L122").Value = 6500.000100000001
$ws.Range("M122").Value = -16667212
$ws.Range("N122").Value = -11400.0001
$ws.Range("H132").Value = 3198
$ws.Range("I132").Value = 2918.0667
$ws.Range("K132").Value = 8754.2001
$ws.Range("M132").Value = -6224.2001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4667.2856
$ws.Range("I61").Value = 5215.0713
$ws.Range("K61").Value = 5215.0713
$ws.Range("M61").Value = -5013.0713
$ws.Range("H113").Value = 4667.2856
$ws.Range("I113").Value = 5215.0713
$ws.Range("K113").Value = 5215.0713
$ws.Range("M113").Value = -3045.0713

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 10015517
$ws.Range("I62").Value = 14304179
$ws.Range("J62").Value = 8637
$ws.Range("K62").Value = 14304179
$ws.Range("L62").Value = 8637
$ws.Range("M62").Value = -14303555
$ws.Range("N62").Value = -9885
$ws.Range("H65").Value = 10015517
$ws.Range("I65").Value = 14304179
$ws.Range("J65").Value = 8637
$ws.Range("K65").Value = 71520895
$ws.Range("L65").Value = 43185
$ws.Range("M65").Value = -71517775
$ws.Range("N65").Value = -49425
$ws.Range("H107").Value = 3473030
$ws.Range("I107").Value = 3968991.5
$ws.Range("J107").Value = 1300
$ws.Range("K107").Value = 11906974.5
$ws.Range("L107").Value = 3900
$ws.Range("M107").Value = -11905054.5
$ws.Range("N107").Value = -7740
$ws.Range("H113").Value = 800
$ws.Range("I113").Value = 683.1539
$ws.Range("J113").Value = 1179.75
$ws.Range("K113").Value = 2049.4617
$ws.Range("L113").Value = 3539.25
$ws.Range("M113").Value = 120.5383000000002
$ws.Range("N113").Value = -7879.25
$ws.Range("H132").Value = 15155521
$ws.Range("J132").Value = 2428.75
$ws.Range("L132").Value = 7286.25
$ws.Range("N132").Value = -12346.25

Write-Host "Applied 255 cell edits"